{"js": "// Map of old (unique) text -> new text, exactly as required by the diff.\nconst replacements = [\n  [\"2025-08-08 Friday\", \"2025-08-09 Saturday\"],\n  [\"589\u00f75=117, 4\", \"514\u00f78=64, 2\"],\n  [\"647\u00f77=92, 3\", \"277\u00f76=46, 1\"],\n  [\"636\u00f76=106, 0\", \"299\u00f72=149, 1\"],\n  [\"310\u00f78=38, 6\", \"522\u00f77=74, 4\"],\n  [\"338\u00f79=37, 5\", \"853\u00f75=170, 3\"],\n  [\"692\u00f76=115, 2\", \"786\u00f72=393, 0\"],\n  [\"423\u00f77=60, 3\", \"355\u00f75=71, 0\"],\n  [\"259\u00f77=37, 0\", \"200\u00f77=28, 4\"],\n  [\"415\u00f74=103, 3\", \"889\u00f75=177, 4\"],\n  [\"369\u00f79=41, 0\", \"285\u00f79=31, 6\"],\n  [\"141\u00f78=17, 5\", \"755\u00f72=377, 1\"],\n  [\"597\u00f77=85, 2\", \"647\u00f79=71, 8\"],\n  [\"587\u00f75=117, 2\", \"546\u00f78=68, 2\"],\n  [\"703\u00f79=78, 1\", \"476\u00f75=95, 1\"],\n  [\"114\u00f73=38, 0\", \"243\u00f77=34, 5\"],\n  [\"299\u00f79=33, 2\", \"736\u00f76=122, 4\"],\n  [\"420\u00f79=46, 6\", \"437\u00f76=72, 5\"],\n  [\"544\u00f79=60, 4\", \"661\u00f73=220, 1\"],\n  [\"432\u00f76=72, 0\", \"486\u00f78=60, 6\"],\n  [\"474\u00f76=79, 0\", \"138\u00f74=34, 2\"],\n  [\"708\u00f77=101, 1\", \"388\u00f78=48, 4\"],\n  [\"914\u00f79=101, 5\", \"648\u00f77=92, 4\"],\n  [\"827\u00f72=413, 1\", \"850\u00f78=106, 2\"],\n  [\"453\u00f73=151, 0\", \"329\u00f79=36, 5\"],\n  [\"137\u00f79=15, 2\", \"677\u00f72=338, 1\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const found = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  found.load(\"items\");\n  await context.sync();\n\n  if (found.items.length === 0) {\n    throw new Error(`Text not found: ${oldText}`);\n  }\n\n  for (const range of found.items) {\n    range.insertText(newText, \"Replace\");\n  }\n}\n\nawait context.sync();\n", "ps1": "# Word COM interop: apply the same 26 text replacements (1 date line +\n# 25 division-problem cells) described by the diff.\n$d = $word.ActiveDocument\n\n$wdReplaceAll = 2\n$wdFindWrapContinue = 1\n\n$pairs = @(\n    @(\"2025-08-08 Friday\", \"2025-08-09 Saturday\"),\n    @(\"589\u00f75=117, 4\", \"514\u00f78=64, 2\"),\n    @(\"647\u00f77=92, 3\", \"277\u00f76=46, 1\"),\n    @(\"636\u00f76=106, 0\", \"299\u00f72=149, 1\"),\n    @(\"310\u00f78=38, 6\", \"522\u00f77=74, 4\"),\n    @(\"338\u00f79=37, 5\", \"853\u00f75=170, 3\"),\n    @(\"692\u00f76=115, 2\", \"786\u00f72=393, 0\"),\n    @(\"423\u00f77=60, 3\", \"355\u00f75=71, 0\"),\n    @(\"259\u00f77=37, 0\", \"200\u00f77=28, 4\"),\n    @(\"415\u00f74=103, 3\", \"889\u00f75=177, 4\"),\n    @(\"369\u00f79=41, 0\", \"285\u00f79=31, 6\"),\n    @(\"141\u00f78=17, 5\", \"755\u00f72=377, 1\"),\n    @(\"597\u00f77=85, 2\", \"647\u00f79=71, 8\"),\n    @(\"587\u00f75=117, 2\", \"546\u00f78=68, 2\"),\n    @(\"703\u00f79=78, 1\", \"476\u00f75=95, 1\"),\n    @(\"114\u00f73=38, 0\", \"243\u00f77=34, 5\"),\n    @(\"299\u00f79=33, 2\", \"736\u00f76=122, 4\"),\n    @(\"420\u00f79=46, 6\", \"437\u00f76=72, 5\"),\n    @(\"544\u00f79=60, 4\", \"661\u00f73=220, 1\"),\n    @(\"432\u00f76=72, 0\", \"486\u00f78=60, 6\"),\n    @(\"474\u00f76=79, 0\", \"138\u00f74=34, 2\"),\n    @(\"708\u00f77=101, 1\", \"388\u00f78=48, 4\"),\n    @(\"914\u00f79=101, 5\", \"648\u00f77=92, 4\"),\n    @(\"827\u00f72=413, 1\", \"850\u00f78=106, 2\"),\n    @(\"453\u00f73=151, 0\", \"329\u00f79=36, 5\"),\n    @(\"137\u00f79=15, 2\", \"677\u00f72=338, 1\")\n)\n\nforeach ($pair in $pairs) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n\n    $range = $d.Content\n    $range.Find.ClearFormatting()\n    $range.Find.Replacement.ClearFormatting()\n    $found = $range.Find.Execute(\n        $oldText,\n        $false,\n        $false,\n        $false,\n        $false,\n        $false,\n        $true,\n        $wdFindWrapContinue,\n        $false,\n        $newText,\n        $wdReplaceAll\n    )\n    if (-not $found) {\n        Write-Output \"NOT FOUND\"\n    }\n}\n"}
